# Refresh the cryptocurrency price / 1h-volume table (scheduled data pull).
# Row order for a few coins also changed upstream (re-ranked), which shows
# up below as whole-row (B/C/D/E) rewrites rather than single-cell edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.238.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.895.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.655'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.46%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.45'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.345'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '50.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.15%  '
$ws.Range("E11").Value = '  +2.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0998'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.170.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.98%  '
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.900.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.226.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0813'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '240.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +31.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.43%  '
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.02%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.941'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +17.92%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0560'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("B34").Value = 'BinanceUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.03'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.90%  '
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.50%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0207'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("E41").Value = '  +14.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.91'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '88.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.338.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +42.17%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.93%  '
$ws.Range("E47").Value = '  -1.94%  '
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.079.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -15.70%  '
